# Auto-generated edit script for CodeSystem-ch-elm-foph-business-rules.xlsx
# Bumps the IG metadata (Version / Title / Date) and rewrites the Concepts
# table with the expanded FOPH-0xx business-rule codes.
$wb = $excel.ActiveWorkbook

# ---- Sheet "Metadata": bump Version / Title / Date ----
$meta = $wb.Worksheets.Item(1)
$meta.Range("B3").Value2 = "1.2.0"
$meta.Range("B5").Value2 = "CH ELM FOPH Business Rules"
$meta.Range("B8").Value2 = "2024-03-28T10:46:20+01:00"

# ---- Sheet "Concepts": rebuild the business-rule concept table ----
$ws = $wb.Worksheets.Item(2)

# Snapshot the two Level values ("1" / "2") into scratch cells *before* any row
# is touched, since nearly every row gets overwritten below. Going through a
# value-only paste (instead of a direct Value2 assignment) keeps these numeric-
# looking strings typed as shared-string text instead of being coerced to numbers.
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("A4").Copy()
$ws.Range("Z2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$levelOne = $ws.Range("Z1")
$levelTwo = $ws.Range("Z2")

# Template rows already carrying the data-row style (border + wrap + vertical-top):
#   A7:D7 -> a fully four-column-populated row
# Used as the source for format-only paste onto every newly-needed row so no new
# cell styles get minted.
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row 2: information / Information
$levelOne.Copy()
$ws.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B2").Value2 = "information"
$ws.Range("C2").Value2 = "Information"
$ws.Range("D2").Value2 = "A purely informational message."

# Row 3: FOPH-000 / Generic Information
$levelTwo.Copy()
$ws.Range("A3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B3").Value2 = "FOPH-000"
$ws.Range("C3").Value2 = "Generic Information"

# Row 4: warning / Warning
$levelOne.Copy()
$ws.Range("A4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B4").Value2 = "warning"
$ws.Range("C4").Value2 = "Warning"
$ws.Range("D4").Value2 = "If the rule is violated, the resource is conformant, but it is not necessarily following best practice."

# Row 5: FOPH-005 / The required anonymization for the reported organism was violated. The following field(s) are affected: %fields%
$levelTwo.Copy()
$ws.Range("A5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B5").Value2 = "FOPH-005"
$ws.Range("C5").Value2 = "The required anonymization for the reported organism was violated. The following field(s) are affected: %fields%"

# Row 6: FOPH-011 / The material is already specified by the leading code. The additional material specified in specimen.type will be ignored.
$levelTwo.Copy()
$ws.Range("A6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B6").Value2 = "FOPH-011"
$ws.Range("C6").Value2 = "The material is already specified by the leading code. The additional material specified in specimen.type will be ignored."

# Row 7: FOPH-010 / Attention, the code %code% (%codeSystem%) expires on %validTo%. Please adjust your systems by this date.
$levelTwo.Copy()
$ws.Range("A7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B7").Value2 = "FOPH-010"
$ws.Range("C7").Value2 = "Attention, the code %code% (%codeSystem%) expires on %validTo%. Please adjust your systems by this date."

# Row 8: FOPH-001 / Generic Warning
$levelTwo.Copy()
$ws.Range("A8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B8").Value2 = "FOPH-001"
$ws.Range("C8").Value2 = "Generic Warning"

# Row 9: FOPH-006 / The following elements for the patient's address are expected: %missingElements%.
$levelTwo.Copy()
$ws.Range("A9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B9").Value2 = "FOPH-006"
$ws.Range("C9").Value2 = "The following elements for the patient's address are expected: %missingElements%."

# Row 10: error / Error
$levelOne.Copy()
$ws.Range("A10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B10").Value2 = "error"
$ws.Range("C10").Value2 = "Error"
$ws.Range("D10").Value2 = "If the rule is violated, the resource is not conformant."

# Row 11: FOPH-009 / The transmitted code %code% (%codeSystem%) is outside the defined validity period %validFrom% - %validTo%.
$levelTwo.Copy()
$ws.Range("A11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B11").Value2 = "FOPH-009"
$ws.Range("C11").Value2 = "The transmitted code %code% (%codeSystem%) is outside the defined validity period %validFrom% - %validTo%."

# Row 12: FOPH-002 / Generic Error
$levelTwo.Copy()
$ws.Range("A12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B12").Value2 = "FOPH-002"
$ws.Range("C12").Value2 = "Generic Error"

# Row 13: FOPH-007 / The specified material is not supported by the provided leading code.
$levelTwo.Copy()
$ws.Range("A13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B13").Value2 = "FOPH-007"
$ws.Range("C13").Value2 = "The specified material is not supported by the provided leading code."

# Row 14: FOPH-012 / The transmitted leading code %code% (%codeSystem%) could not be found in the current value set.
$levelTwo.Copy()
$ws.Range("A14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B14").Value2 = "FOPH-012"
$ws.Range("C14").Value2 = "The transmitted leading code %code% (%codeSystem%) could not be found in the current value set."

# Row 15: FOPH-008 / The specified organism is not supported by the provided leading code.
$levelTwo.Copy()
$ws.Range("A15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B15").Value2 = "FOPH-008"
$ws.Range("C15").Value2 = "The specified organism is not supported by the provided leading code."

# Row 16: FOPH-004 / The provided laboratory identification "%identifier%" is either unknown or corresponds to a laboratory for which your account does not have reporting permissions. Please verify your laboratory identification. If it's correct, complete the necessary onboarding process before submitting data on behalf of this lab.
$levelTwo.Copy()
$ws.Range("A16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B16").Value2 = "FOPH-004"
$ws.Range("C16").Value2 = "The provided laboratory identification ""%identifier%"" is either unknown or corresponds to a laboratory for which your account does not have reporting permissions. Please verify your laboratory identification. If it's correct, complete the necessary onboarding process before submitting data on behalf of this lab."

# Row 17: FOPH-013 / The transmitted code %code% (%codeSystem%) is outside the defined validity period.
$levelTwo.Copy()
$ws.Range("A17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("B17").Value2 = "FOPH-013"
$ws.Range("C17").Value2 = "The transmitted code %code% (%codeSystem%) is outside the defined validity period."

# Drop the scratch cells used to stage the Level values - they must not survive
# in the saved sheet / dimension.
$ws.Range("Z1:Z2").Delete()
$excel.CutCopyMode = $false

Write-Host "done"
